# repull data, push all data, mean calculation
# Updates column F (dSF) values on Sheet1 for the rows whose source data
# was repulled, per the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -3
    4  = 3
    11 = -3
    12 = -2
    13 = -4
    16 = 0
    22 = 1
    23 = -1
    25 = 3
    26 = 4
    30 = 0
    33 = 3
    40 = 0
    41 = -1
    45 = 0
    59 = 0
    60 = 0
    62 = 1
    64 = -2
    67 = 3
    70 = 1
    73 = -6
    77 = -1
    78 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
